# Update workbook "horarios-141-completo.xlsx" with the 31/12/2025 13:28 scrape
# for routes LP1912 (sheet "LP1912"), LP1912-215 (sheet "LP1912-215") and
# 6203-6173 (sheet "6203-6173").
#
# LP1912         -> 14 new rows of arrival data (rows 965-978)
# LP1912-215     -> no new arrival rows this run, only the refresh timestamp
# 6203-6173      -> 2 new rows of arrival data (rows 119-120)

$wb = $excel.ActiveWorkbook

$updateStamp = "Última actualización: 31/12/2025 13:28:40"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912   (columns: A=title/blank, B=Hora_Scrap, C=Hora_Llegada,
#                     D=Línea, E=Minutos, F=Parada, G=Fecha)
# ---------------------------------------------------------------------------
$wsLP = $wb.Worksheets.Item("LP1912")

$wsLP.Cells.Item(2, 1).Value = $updateStamp
$wsLP.Cells.Item(3, 1).Value = "Total filas: 977"

$lp1912Rows = @(
  @(965, "13:28:29", "13:31", "16_P MOR-SANTA ANA",    3, "LP1912", "31/12/2025"),
  @(966, "13:28:29", "13:34", "23_HERNANDEZ",           6, "LP1912", "31/12/2025"),
  @(967, "13:28:29", "13:41", "16_SANTA ANA",          13, "LP1912", "31/12/2025"),
  @(968, "13:28:29", "13:51", "15_ABASTO",             23, "LP1912", "31/12/2025"),
  @(969, "13:28:29", "13:51", "16_SANTA ANA",          23, "LP1912", "31/12/2025"),
  @(970, "13:28:29", "14:01", "16_SANTA ANA",          33, "LP1912", "31/12/2025"),
  @(971, "13:28:29", "14:01", "17_ROMERO",             33, "LP1912", "31/12/2025"),
  @(972, "13:28:29", "14:03", "23_HERNANDEZ",          35, "LP1912", "31/12/2025"),
  @(973, "13:28:29", "14:11", "15_ABASTO",             43, "LP1912", "31/12/2025"),
  @(974, "13:28:29", "14:24", "11_ETCHEVERRY",         56, "LP1912", "31/12/2025"),
  @(975, "13:28:29", "14:37", "16_P MOR-SANTA ANA",    69, "LP1912", "31/12/2025"),
  @(976, "13:28:29", "14:39", "23_HERNANDEZ",          71, "LP1912", "31/12/2025"),
  @(977, "13:28:29", "14:40", "17X38_ROMERO",          72, "LP1912", "31/12/2025"),
  @(978, "13:28:29", "15:04", "14_ABASTO",             96, "LP1912", "31/12/2025")
)

foreach ($row in $lp1912Rows) {
  $r = $row[0]
  $wsLP.Cells.Item($r, 1).Value = ""
  $wsLP.Cells.Item($r, 2).Value = $row[1]
  $wsLP.Cells.Item($r, 3).Value = $row[2]
  $wsLP.Cells.Item($r, 4).Value = $row[3]
  $wsLP.Cells.Item($r, 5).Value = $row[4]
  $wsLP.Cells.Item($r, 6).Value = $row[5]
  $wsLP.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215   (columns: A=title/blank, B=Fecha, C=Hora_Scrap,
#                         D=Hora_Llegada, E=Línea, F=Minutos, G=Parada)
# Only the refresh timestamp changes this round - no new arrivals scraped.
# ---------------------------------------------------------------------------
$wsLP215 = $wb.Worksheets.Item("LP1912-215")
$wsLP215.Cells.Item(2, 1).Value = $updateStamp

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173   (columns: A=title/blank, B=Fecha, C=Hora_Scrap,
#                        D=Hora_Llegada, E=Línea, F=Minutos, G=Parada)
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("6203-6173")

$wsB.Cells.Item(2, 1).Value = $updateStamp
$wsB.Cells.Item(3, 1).Value = "Total filas: 119"

$busRows = @(
  @(119, "31/12/2025", "13:28:34", "13:54", "215C_LA PLATA", 26, "L6203"),
  @(120, "31/12/2025", "13:28:34", "14:34", "215C_LA PLATA", 66, "L6203")
)

foreach ($row in $busRows) {
  $r = $row[0]
  $wsB.Cells.Item($r, 1).Value = ""
  $wsB.Cells.Item($r, 2).Value = $row[1]
  $wsB.Cells.Item($r, 3).Value = $row[2]
  $wsB.Cells.Item($r, 4).Value = $row[3]
  $wsB.Cells.Item($r, 5).Value = $row[4]
  $wsB.Cells.Item($r, 6).Value = $row[5]
  $wsB.Cells.Item($r, 7).Value = $row[6]
}

Write-Output "edit.ps1 completed"
